$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 984.0952
$ws.Range("J17").Value = 988.5854
$ws.Range("L17").Value = 2965.7562
$ws.Range("N17").Value = -3301.7562
$ws.Range("H32").Value = 6486.75
$ws.Range("I32").Value = 4833.3335
$ws.Range("K32").Value = 4833.3335
$ws.Range("M32").Value = -4507.3335
$ws.Range("J62").Value = 88877.14
$ws.Range("L62").Value = 88877.14
$ws.Range("N62").Value = -90125.14
$ws.Range("H64").Value = 6994.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 6994.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6994.5
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -7490.5
$ws.Range("J65").Value = 88877.14
$ws.Range("L65").Value = 444385.7
$ws.Range("N65").Value = -450625.7
$ws.Range("H67").Value = 6994.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 6994.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6994.5
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -8710.5
$ws.Range("H70").Value = 32409564
$ws.Range("I70").Value = 50002310
$ws.Range("J70").Value = 25643122
$ws.Range("K70").Value = 150006930
$ws.Range("L70").Value = 76929366
$ws.Range("M70").Value = -150006660
$ws.Range("N70").Value = -76929906
$ws.Range("H73").Value = 32409564
$ws.Range("I73").Value = 50002310
$ws.Range("J73").Value = 25643122
$ws.Range("K73").Value = 150006930
$ws.Range("L73").Value = 76929366
$ws.Range("M73").Value = -150005994
$ws.Range("N73").Value = -76931238
$ws.Range("H92").Value = 1336.4375
$ws.Range("I92").Value = 773.3
$ws.Range("J92").Value = 2275
$ws.Range("K92").Value = 773.3
$ws.Range("L92").Value = 2275
$ws.Range("M92").Value = 474.7
$ws.Range("N92").Value = -4771
$ws.Range("H132").Value = 1558.3334
$ws.Range("I132").Value = 1558.3334
$ws.Range("K132").Value = 4675.0002
$ws.Range("M132").Value = -2145.0002
$ws.Range("H137").Value = 2631.6052
$ws.Range("I137").Value = 2345.8696
$ws.Range("J137").Value = 3069.7334
$ws.Range("K137").Value = 7037.6088
$ws.Range("L137").Value = 9209.200199999999
$ws.Range("M137").Value = -4487.6088
$ws.Range("N137").Value = -14309.2002
$ws.Range("H138").Value = 3803.2886
$ws.Range("I138").Value = 1181.9656
$ws.Range("J138").Value = 7108.4346
$ws.Range("K138").Value = 3545.8968
$ws.Range("L138").Value = 21325.3038
$ws.Range("M138").Value = 1594.1032
$ws.Range("N138").Value = -31605.3038

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6264.8613
$ws.Range("I61").Value = 2617.6365
$ws.Range("J61").Value = 11996.214
$ws.Range("K61").Value = 2617.6365
$ws.Range("L61").Value = 11996.214
$ws.Range("M61").Value = -2405.6365
$ws.Range("N61").Value = -12420.214
$ws.Range("H63").Value = 2287
$ws.Range("J63").Value = 2276.5
$ws.Range("L63").Value = 2276.5
$ws.Range("N63").Value = -3648.5
$ws.Range("H66").Value = 2287
$ws.Range("J66").Value = 2276.5
$ws.Range("L66").Value = 11382.5
$ws.Range("N66").Value = -18246.5
$ws.Range("H74").Value = 16929.256
$ws.Range("I74").Value = 30400.285
$ws.Range("J74").Value = 4070.5454
$ws.Range("K74").Value = 30400.285
$ws.Range("L74").Value = 4070.5454
$ws.Range("M74").Value = -29526.285
$ws.Range("N74").Value = -5818.5454
$ws.Range("H77").Value = 16929.256
$ws.Range("I77").Value = 30400.285
$ws.Range("J77").Value = 4070.5454
$ws.Range("K77").Value = 152001.425
$ws.Range("L77").Value = 20352.727
$ws.Range("M77").Value = -147633.425
$ws.Range("N77").Value = -29088.727
$ws.Range("H102").Value = 899.9048
$ws.Range("I102").Value = 839.3684
$ws.Range("K102").Value = 839.3684
$ws.Range("M102").Value = 782.6316
$ws.Range("H136").Value = 6264.8613
$ws.Range("I136").Value = 2617.6365
$ws.Range("J136").Value = 11996.214
$ws.Range("K136").Value = 7852.9095
$ws.Range("L136").Value = 35988.642
$ws.Range("M136").Value = -5302.9095
$ws.Range("N136").Value = -41088.642

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2600542.2
$ws.Range("I99").Value = 2981.5518
$ws.Range("K99").Value = 2981.5518
$ws.Range("M99").Value = -1483.5518
$ws.Range("H105").Value = 48062.06
$ws.Range("J105").Value = 3927.5715
$ws.Range("L105").Value = 3927.5715
$ws.Range("N105").Value = -7421.5715

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6704
$ws.Range("I31").Value = 2410.8823
$ws.Range("K31").Value = 2410.8823
$ws.Range("M31").Value = -2115.8823
$ws.Range("H34").Value = 6704
$ws.Range("I34").Value = 2410.8823
$ws.Range("K34").Value = 2410.8823
$ws.Range("M34").Value = -2208.8823
$ws.Range("H58").Value = 14712677
$ws.Range("I58").Value = 35715884
$ws.Range("K58").Value = 35715884
$ws.Range("M58").Value = -35715681
$ws.Range("H109").Value = 50285
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 50285
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 50285
$ws.Range("M109").Value = ""
$ws.Range("N109").Value = -52365
$ws.Range("H132").Value = 6080.8857
$ws.Range("I132").Value = 2481.9285
$ws.Range("J132").Value = 8480.190000000001
$ws.Range("K132").Value = 7445.7855
$ws.Range("L132").Value = 25440.57
$ws.Range("M132").Value = -4915.7855
$ws.Range("N132").Value = -30500.57
$ws.Range("H134").Value = 5993.282
$ws.Range("I134").Value = 1652.4375
$ws.Range("K134").Value = 4957.3125
$ws.Range("M134").Value = -2422.3125
$ws.Range("H136").Value = 14712677
$ws.Range("I136").Value = 35715884
$ws.Range("K136").Value = 107147652
$ws.Range("M136").Value = -107145102

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 623.73334
$ws.Range("I114").Value = 401.4
$ws.Range("J114").Value = 734.9
$ws.Range("K114").Value = 1204.2
$ws.Range("L114").Value = 2204.7
$ws.Range("M114").Value = 2049.8
$ws.Range("N114").Value = -8712.700000000001
$ws.Range("H137").Value = 54416.844
$ws.Range("I137").Value = 1437.8
$ws.Range("J137").Value = 113282.445
$ws.Range("K137").Value = 4313.4
$ws.Range("L137").Value = 339847.335
$ws.Range("M137").Value = 786.6000000000004
$ws.Range("N137").Value = -350047.335

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6378.8125
$ws.Range("I132").Value = 1912.6666
$ws.Range("K132").Value = 5737.9998
$ws.Range("M132").Value = -3207.9998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6270.0713
$ws.Range("I122").Value = 4222
$ws.Range("K122").Value = 12666
$ws.Range("M122").Value = -10216
$ws.Range("H136").Value = 10003.893
$ws.Range("I136").Value = 2139
$ws.Range("J136").Value = 13149.85
$ws.Range("K136").Value = 6417
$ws.Range("L136").Value = 39449.55
$ws.Range("M136").Value = -3867
$ws.Range("N136").Value = -44549.55

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8994
$ws.Range("I62").Value = 8994
$ws.Range("K62").Value = 8994
$ws.Range("M62").Value = -8370
$ws.Range("H65").Value = 8994
$ws.Range("I65").Value = 8994
$ws.Range("K65").Value = 44970
$ws.Range("M65").Value = -41850
$ws.Range("H81").Value = 33346580
$ws.Range("I81").Value = 4999
$ws.Range("J81").Value = 66688160
$ws.Range("K81").Value = 9998
$ws.Range("L81").Value = 133376320
$ws.Range("M81").Value = -8937
$ws.Range("N81").Value = -133378442
$ws.Range("H84").Value = 33346580
$ws.Range("I84").Value = 4999
$ws.Range("J84").Value = 66688160
$ws.Range("K84").Value = 49990
$ws.Range("L84").Value = 666881600
$ws.Range("M84").Value = -44686
$ws.Range("N84").Value = -666892208
$ws.Range("H132").Value = 13900781
$ws.Range("H136").Value = 27059258
$ws.Range("I136").Value = 76924240
$ws.Range("J136").Value = 49059.293
$ws.Range("K136").Value = 230772720
$ws.Range("L136").Value = 147177.879
$ws.Range("M136").Value = -230770170
$ws.Range("N136").Value = -152277.879
